$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.104.98"
$ws.Range("E2").Value = "  +3.39%  "

$ws.Range("D3").Value = "3.482.75"
$ws.Range("E3").Value = "  +3.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.71"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.91"
$ws.Range("E6").Value = "  +17.66%  "

$ws.Range("D7").Value = "3.473.49"
$ws.Range("E7").Value = "  +3.24%  "

$ws.Range("E8").Value = "  +2.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.701"
$ws.Range("E10").Value = "  +8.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.132"
$ws.Range("E11").Value = "  +32.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.58"
$ws.Range("E12").Value = "  +9.80%  "

$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("D14").Value = "4.021.17"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.84"
$ws.Range("E15").Value = "  +4.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.21"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "3.450.33"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("D18").Value = "63.011.10"
$ws.Range("E18").Value = "  +3.51%  "

$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.85"
$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000145"
$ws.Range("E21").Value = "  +29.76%  "

$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.98"
$ws.Range("E23").Value = "  +10.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.18"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "312.99"
$ws.Range("E25").Value = "  +3.00%  "

$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.62"
$ws.Range("E27").Value = "  +5.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("E28").Value = "  +2.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.178"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("E30").Value = "  -2.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.57"
$ws.Range("E31").Value = "  -3.56%  "

$ws.Range("E32").Value = "  +2.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "44.30"
$ws.Range("E33").Value = "  +11.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.90"
$ws.Range("E34").Value = "  +3.55%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0496"
$ws.Range("E37").Value = "  -2.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.72"
$ws.Range("E38").Value = "  +0.53%  "

$ws.Range("E39").Value = "  +5.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.05"
$ws.Range("E41").Value = "  -3.16%  "

$ws.Range("E42").Value = "  +2.40%  "

$ws.Range("E43").Value = "  +3.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.49"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.54"
$ws.Range("E45").Value = "  +3.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.289"
$ws.Range("E46").Value = "  -5.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.99"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.22"
$ws.Range("E49").Value = "  -1.86%  "

$ws.Range("D50").Value = "3.823.36"
$ws.Range("E50").Value = "  +3.70%  "

$ws.Range("D51").Value = "2.189.41"
$ws.Range("E51").Value = "  -0.10%  "
